# Weekly update: insert two new price rows for the most recent week, pushing
# all existing data rows down by two (rows 574-686 -> 576-688).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 574 (entire rows), shifting rows 574:686 down to 576:688.
$ws.Range("A574:A575").EntireRow.Insert()

# Populate the first new row (574) with the new weekly data.
$ws.Range("A574").Value = 6
$ws.Range("B574").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C574").Value = "Metropolitana"
$ws.Range("D574").Value = 44637
$ws.Range("E574").Value = 13
$ws.Range("F574").Value = 100112023
$ws.Range("G574").Value = "Brócoli"
$ws.Range("H574").Value = "Sin especificar"
$ws.Range("I574").Value = "Primera"
$ws.Range("J574").Value = 13400
$ws.Range("K574").Value = 850
$ws.Range("L574").Value = 1000
$ws.Range("M574").Value = 935
$ws.Range("N574").Value = "$/unidad"
$ws.Range("O574").Value = "Región Metropolitana"
$ws.Range("P574").Value = 935
$ws.Range("Q574").Value = 1
$ws.Range("R574").Value = "Hortaliza"

# Populate the second new row (575) with the new weekly data.
$ws.Range("A575").Value = 6
$ws.Range("B575").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C575").Value = "Metropolitana"
$ws.Range("D575").Value = 44637
$ws.Range("E575").Value = 13
$ws.Range("F575").Value = 100112023
$ws.Range("G575").Value = "Brócoli"
$ws.Range("H575").Value = "Sin especificar"
$ws.Range("I575").Value = "Segunda"
$ws.Range("J575").Value = 5700
$ws.Range("K575").Value = 600
$ws.Range("L575").Value = 700
$ws.Range("M575").Value = 642
$ws.Range("N575").Value = "$/unidad"
$ws.Range("O575").Value = "Región Metropolitana"
$ws.Range("P575").Value = 642
$ws.Range("Q575").Value = 1
$ws.Range("R575").Value = "Hortaliza"
